# Applies the updated Leve profit-calculation figures captured by the scheduled
# Sheets runner (currentAveragePrice/-NQ/-HQ, LevePriceNQ/HQ and profit columns).
$wb = $excel.ActiveWorkbook

$ws_ALC = $wb.Worksheets.Item("ALC")
# Row 62
$ws_ALC.Range("H62").Value = 7482.2383
$ws_ALC.Range("I62").Value = 6268.375
$ws_ALC.Range("K62").Value = 6268.375
$ws_ALC.Range("M62").Value = -5644.375
# Row 65
$ws_ALC.Range("H65").Value = 7482.2383
$ws_ALC.Range("I65").Value = 6268.375
$ws_ALC.Range("K65").Value = 31341.875
$ws_ALC.Range("M65").Value = -28221.875
# Row 69
$ws_ALC.Range("H69").Value = 7622.615
$ws_ALC.Range("J69").Value = 8546.5
$ws_ALC.Range("L69").Value = 25639.5
$ws_ALC.Range("N69").Value = -27387.5
# Row 72
$ws_ALC.Range("H72").Value = 7622.615
$ws_ALC.Range("J72").Value = 8546.5
$ws_ALC.Range("L72").Value = 76918.5
$ws_ALC.Range("N72").Value = -85654.5
# Row 95
$ws_ALC.Range("H95").Value = 76500
$ws_ALC.Range("J95").Value = 76500
$ws_ALC.Range("L95").Value = 76500
$ws_ALC.Range("N95").Value = -81992
# Row 112
$ws_ALC.Range("H112").Value = 1710.8077
$ws_ALC.Range("J112").Value = 1719.24
$ws_ALC.Range("L112").Value = 5157.72
$ws_ALC.Range("N112").Value = -7373.72
# Row 128
$ws_ALC.Range("H128").Value = 36294.43
$ws_ALC.Range("J128").Value = 32978
$ws_ALC.Range("L128").Value = 32978
$ws_ALC.Range("N128").Value = -42938
# Row 137
$ws_ALC.Range("H137").Value = 13160687
$ws_ALC.Range("J137").Value = 2850.484
$ws_ALC.Range("L137").Value = 8551.451999999999
$ws_ALC.Range("N137").Value = -13651.452
# Row 138
$ws_ALC.Range("H138").Value = 3512.9285
$ws_ALC.Range("I138").Value = 1871.4706
$ws_ALC.Range("J138").Value = 4039.434
$ws_ALC.Range("K138").Value = 5614.4118
$ws_ALC.Range("L138").Value = 12118.302
$ws_ALC.Range("M138").Value = -474.4117999999999
$ws_ALC.Range("N138").Value = -22398.302

$ws_ARM = $wb.Worksheets.Item("ARM")
# Row 33
$ws_ARM.Range("H33").Value = 13999.571
$ws_ARM.Range("I33").Value = 5999
$ws_ARM.Range("K33").Value = 5999
$ws_ARM.Range("M33").Value = -5670
# Row 36
$ws_ARM.Range("H36").Value = 11947.909
$ws_ARM.Range("I36").Value = 5237.8335
$ws_ARM.Range("K36").Value = 5237.8335
$ws_ARM.Range("M36").Value = -4891.8335
# Row 97
$ws_ARM.Range("H97").Value = 890.5625
$ws_ARM.Range("I97").Value = 907
$ws_ARM.Range("K97").Value = 907
$ws_ARM.Range("M97").Value = -411

$ws_CRP = $wb.Worksheets.Item("CRP")
# Row 16
$ws_CRP.Range("H16").Value = 1750.8889
$ws_CRP.Range("I16").Value = 879.7778
$ws_CRP.Range("J16").Value = 2622
$ws_CRP.Range("K16").Value = 879.7778
$ws_CRP.Range("L16").Value = 2622
$ws_CRP.Range("M16").Value = -592.7778
$ws_CRP.Range("N16").Value = -3196
# Row 31
$ws_CRP.Range("H31").Value = 39340.758
$ws_CRP.Range("I31").Value = 2343.6428
$ws_CRP.Range("K31").Value = 2343.6428
$ws_CRP.Range("M31").Value = -2048.6428
# Row 34
$ws_CRP.Range("H34").Value = 39340.758
$ws_CRP.Range("I34").Value = 2343.6428
$ws_CRP.Range("K34").Value = 2343.6428
$ws_CRP.Range("M34").Value = -2141.6428
# Row 105
$ws_CRP.Range("H105").Value = 2408.2666
$ws_CRP.Range("I105").Value = 912.1818
$ws_CRP.Range("J105").Value = 6522.5
$ws_CRP.Range("K105").Value = 912.1818
$ws_CRP.Range("L105").Value = 6522.5
$ws_CRP.Range("M105").Value = 834.8182
$ws_CRP.Range("N105").Value = -10016.5
# Row 113
$ws_CRP.Range("H113").Value = 1750.8889
$ws_CRP.Range("I113").Value = 879.7778
$ws_CRP.Range("J113").Value = 2622
$ws_CRP.Range("K113").Value = 879.7778
$ws_CRP.Range("L113").Value = 2622
$ws_CRP.Range("M113").Value = 1290.2222
$ws_CRP.Range("N113").Value = -6962
# Row 132
$ws_CRP.Range("H132").Value = 2956.9019
$ws_CRP.Range("I132").Value = 2500.125
$ws_CRP.Range("K132").Value = 7500.375
$ws_CRP.Range("M132").Value = -4970.375

$ws_CUL = $wb.Worksheets.Item("CUL")
# Row 7
$ws_CUL.Range("H7").Value = 83720
$ws_CUL.Range("I7").Value = 660
$ws_CUL.Range("J7").Value = 166780
$ws_CUL.Range("K7").Value = 1980
$ws_CUL.Range("L7").Value = 500340
$ws_CUL.Range("M7").Value = -1868
$ws_CUL.Range("N7").Value = -500564
# Row 113
$ws_CUL.Range("H113").Value = 100001810
$ws_CUL.Range("I113").Value = 1499.5
$ws_CUL.Range("K113").Value = 4498.5
$ws_CUL.Range("M113").Value = -2328.5

$ws_GSM = $wb.Worksheets.Item("GSM")
# Row 97
$ws_GSM.Range("H97").Value = 3747.3333
$ws_GSM.Range("I97").Value = 3495
$ws_GSM.Range("J97").Value = 3873.5
$ws_GSM.Range("K97").Value = 3495
$ws_GSM.Range("L97").Value = 3873.5
$ws_GSM.Range("M97").Value = -2999
$ws_GSM.Range("N97").Value = -4865.5

$ws_LTW = $wb.Worksheets.Item("LTW")
# Row 16
$ws_LTW.Range("H16").Value = 1340.9546
$ws_LTW.Range("I16").Value = 475.1
$ws_LTW.Range("J16").Value = 9999.5
$ws_LTW.Range("K16").Value = 475.1
$ws_LTW.Range("L16").Value = 9999.5
$ws_LTW.Range("M16").Value = -305.1
$ws_LTW.Range("N16").Value = -10339.5
# Row 61
$ws_LTW.Range("H61").Value = 8609.177
$ws_LTW.Range("I61").Value = 7946.1665
$ws_LTW.Range("J61").Value = 10200.4
$ws_LTW.Range("K61").Value = 7946.1665
$ws_LTW.Range("L61").Value = 10200.4
$ws_LTW.Range("M61").Value = -7744.1665
$ws_LTW.Range("N61").Value = -10604.4
# Row 113
$ws_LTW.Range("H113").Value = 8609.177
$ws_LTW.Range("I113").Value = 7946.1665
$ws_LTW.Range("J113").Value = 10200.4
$ws_LTW.Range("K113").Value = 7946.1665
$ws_LTW.Range("L113").Value = 10200.4
$ws_LTW.Range("M113").Value = -5776.1665
$ws_LTW.Range("N113").Value = -14540.4
# Row 123
$ws_LTW.Range("H123").Value = 119999
$ws_LTW.Range("J123").Value = 119999
$ws_LTW.Range("L123").Value = 119999
$ws_LTW.Range("N123").Value = -129799
# Row 132
$ws_LTW.Range("H132").Value = 8071.933
$ws_LTW.Range("I132").Value = 13001.333
$ws_LTW.Range("J132").Value = 6839.5835
$ws_LTW.Range("K132").Value = 39003.999
$ws_LTW.Range("L132").Value = 20518.7505
$ws_LTW.Range("M132").Value = -36473.999
$ws_LTW.Range("N132").Value = -25578.7505
# Row 136
$ws_LTW.Range("H136").Value = 3280.0442
$ws_LTW.Range("I136").Value = 2603.641
$ws_LTW.Range("J136").Value = 4189.6895
$ws_LTW.Range("K136").Value = 7810.923000000001
$ws_LTW.Range("L136").Value = 12569.0685
$ws_LTW.Range("M136").Value = -5260.923000000001
$ws_LTW.Range("N136").Value = -17669.0685

$ws_WVR = $wb.Worksheets.Item("WVR")
# Row 125
$ws_WVR.Range("H125").Value = 95000
$ws_WVR.Range("J125").Value = 95000
$ws_WVR.Range("L125").Value = 95000
$ws_WVR.Range("N125").Value = -104840
# Row 127
$ws_WVR.Range("H127").Value = 0
$ws_WVR.Range("J127").Value = 0
$ws_WVR.Range("L127").Value = 0
$ws_WVR.Range("N127").ClearContents()
# Row 132
$ws_WVR.Range("H132").Value = 2790.6562
$ws_WVR.Range("I132").Value = 2295.8696
$ws_WVR.Range("J132").Value = 4055.111
$ws_WVR.Range("K132").Value = 6887.6088
$ws_WVR.Range("L132").Value = 12165.333
$ws_WVR.Range("M132").Value = -4357.6088
$ws_WVR.Range("N132").Value = -17225.333
# Row 136
$ws_WVR.Range("H136").Value = 3076.9412
$ws_WVR.Range("I136").Value = 1056.8182
$ws_WVR.Range("K136").Value = 3170.4546
$ws_WVR.Range("M136").Value = -620.4546
